$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnas"
$ws.Range("C2").Value = "Vipr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 190.664594
$ws.Range("H2").Value = 571.993782
$ws.Range("I2").Value = 0.2001939625490346
$ws.Range("J2").Value = 0.2001939625490346
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.092341333333333
$ws.Range("N2").Value = 3.277024
$ws.Range("O2").Value = 0.1566881615238948
$ws.Range("P2").Value = 0.1566881615238948
$ws.Range("Q2").Value = 208.2708168294187
$ws.Range("R2").Value = 1874.437351464768
$ws.Range("S2").Value = 0.03136802393999168
$ws.Range("T2").Value = 0.03136802393999168

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnas"
$ws.Range("C3").Value = "Vipr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 190.664594
$ws.Range("H3").Value = 571.993782
$ws.Range("I3").Value = 0.2001939625490346
$ws.Range("J3").Value = 0.2001939625490346
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.146504
$ws.Range("N3").Value = 0.439512
$ws.Range("O3").Value = 0.02101489865429428
$ws.Range("P3").Value = 0.02101489865429428
$ws.Range("Q3").Value = 27.933125679376
$ws.Range("R3").Value = 251.398131114384
$ws.Range("S3").Value = 0.004207055834169546
$ws.Range("T3").Value = 0.004207055834169546

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnas"
$ws.Range("C4").Value = "Vipr1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 190.664594
$ws.Range("H4").Value = 571.993782
$ws.Range("I4").Value = 0.2001939625490346
$ws.Range("J4").Value = 0.2001939625490346
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.199477333333333
$ws.Range("N4").Value = 9.598431999999999
$ws.Range("O4").Value = 0.4589409975612385
$ws.Range("P4").Value = 0.4589409975612386
$ws.Range("Q4").Value = 610.0270467722025
$ws.Range("R4").Value = 5490.243420949823
$ws.Range("S4").Value = 0.09187721687799118
$ws.Range("T4").Value = 0.09187721687799119

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Gnas"
$ws.Range("C5").Value = "Vipr1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 190.664594
$ws.Range("H5").Value = 571.993782
$ws.Range("I5").Value = 0.2001939625490346
$ws.Range("J5").Value = 0.2001939625490346
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.533112333333334
$ws.Range("N5").Value = 7.599337000000001
$ws.Range("O5").Value = 0.3633559422605724
$ws.Range("P5").Value = 0.3633559422605724
$ws.Range("Q5").Value = 482.9748345913928
$ws.Range("R5").Value = 4346.773511322534
$ws.Range("S5").Value = 0.07274166589688223
$ws.Range("T5").Value = 0.07274166589688223

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnas"
$ws.Range("C6").Value = "Vipr1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 531.1103823333333
$ws.Range("H6").Value = 1593.331147
$ws.Range("I6").Value = 0.5576551459273178
$ws.Range("J6").Value = 0.5576551459273177
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.092341333333333
$ws.Range("N6").Value = 3.277024
$ws.Range("O6").Value = 0.1566881615238948
$ws.Range("P6").Value = 0.1566881615238948
$ws.Range("Q6").Value = 580.1538231851698
$ws.Range("R6").Value = 5221.384408666528
$ws.Range("S6").Value = 0.08737795957969069
$ws.Range("T6").Value = 0.08737795957969066

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnas"
$ws.Range("C7").Value = "Vipr1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 531.1103823333333
$ws.Range("H7").Value = 1593.331147
$ws.Range("I7").Value = 0.5576551459273178
$ws.Range("J7").Value = 0.5576551459273177
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.146504
$ws.Range("N7").Value = 0.439512
$ws.Range("O7").Value = 0.02101489865429428
$ws.Range("P7").Value = 0.02101489865429428
$ws.Range("Q7").Value = 77.80979545336265
$ws.Range("R7").Value = 700.288159080264
$ws.Range("S7").Value = 0.01171906637570827
$ws.Range("T7").Value = 0.01171906637570827

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Gnas"
$ws.Range("C8").Value = "Vipr1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 531.1103823333333
$ws.Range("H8").Value = 1593.331147
$ws.Range("I8").Value = 0.5576551459273178
$ws.Range("J8").Value = 0.5576551459273177
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.199477333333333
$ws.Range("N8").Value = 9.598431999999999
$ws.Range("O8").Value = 0.4589409975612385
$ws.Range("P8").Value = 0.4589409975612386
$ws.Range("Q8").Value = 1699.2756297735
$ws.Range("R8").Value = 15293.4806679615
$ws.Range("S8").Value = 0.2559308089670413
$ws.Range("T8").Value = 0.2559308089670412

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Gnas"
$ws.Range("C9").Value = "Vipr1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 531.1103823333333
$ws.Range("H9").Value = 1593.331147
$ws.Range("I9").Value = 0.5576551459273178
$ws.Range("J9").Value = 0.5576551459273177
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.533112333333334
$ws.Range("N9").Value = 7.599337000000001
$ws.Range("O9").Value = 0.3633559422605724
$ws.Range("P9").Value = 0.3633559422605724
$ws.Range("Q9").Value = 1345.362259849949
$ws.Range("R9").Value = 12108.26033864954
$ws.Range("S9").Value = 0.2026273110048776
$ws.Range("T9").Value = 0.2026273110048775

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Gnas"
$ws.Range("C10").Value = "Vipr1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 79.06597733333334
$ws.Range("H10").Value = 237.197932
$ws.Range("I10").Value = 0.08301767503395074
$ws.Range("J10").Value = 0.08301767503395074
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.092341333333333
$ws.Range("N10").Value = 3.277024
$ws.Range("O10").Value = 0.1566881615238948
$ws.Range("P10").Value = 0.1566881615238948
$ws.Range("Q10").Value = 86.36703510159646
$ws.Range("R10").Value = 777.303315914368
$ws.Range("S10").Value = 0.01300788687505788
$ws.Range("T10").Value = 0.01300788687505788

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Gnas"
$ws.Range("C11").Value = "Vipr1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 79.06597733333334
$ws.Range("H11").Value = 237.197932
$ws.Range("I11").Value = 0.08301767503395074
$ws.Range("J11").Value = 0.08301767503395074
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.146504
$ws.Range("N11").Value = 0.439512
$ws.Range("O11").Value = 0.02101489865429428
$ws.Range("P11").Value = 0.02101489865429428
$ws.Range("Q11").Value = 11.58348194324267
$ws.Range("R11").Value = 104.251337489184
$ws.Range("S11").Value = 0.001744608027353611
$ws.Range("T11").Value = 0.001744608027353611

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Gnas"
$ws.Range("C12").Value = "Vipr1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 79.06597733333334
$ws.Range("H12").Value = 237.197932
$ws.Range("I12").Value = 0.08301767503395074
$ws.Range("J12").Value = 0.08301767503395074
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.199477333333333
$ws.Range("N12").Value = 9.598431999999999
$ws.Range("O12").Value = 0.4589409975612385
$ws.Range("P12").Value = 0.4589409975612386
$ws.Range("Q12").Value = 252.9698023158471
$ws.Range("R12").Value = 2276.728220842624
$ws.Range("S12").Value = 0.03810021459529608
$ws.Range("T12").Value = 0.03810021459529608

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Gnas"
$ws.Range("C13").Value = "Vipr1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 79.06597733333334
$ws.Range("H13").Value = 237.197932
$ws.Range("I13").Value = 0.08301767503395074
$ws.Range("J13").Value = 0.08301767503395074
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.533112333333334
$ws.Range("N13").Value = 7.599337000000001
$ws.Range("O13").Value = 0.3633559422605724
$ws.Range("P13").Value = 0.3633559422605724
$ws.Range("Q13").Value = 200.2830023301205
$ws.Range("R13").Value = 1802.547020971084
$ws.Range("S13").Value = 0.03016496553624317
$ws.Range("T13").Value = 0.03016496553624317

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Gnas"
$ws.Range("C14").Value = "Vipr1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 151.5583673333333
$ws.Range("H14").Value = 454.675102
$ws.Range("I14").Value = 0.1591332164896969
$ws.Range("J14").Value = 0.1591332164896969
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.092341333333333
$ws.Range("N14").Value = 3.277024
$ws.Range("O14").Value = 0.1566881615238948
$ws.Range("P14").Value = 0.1566881615238948
$ws.Range("Q14").Value = 165.5534690507164
$ws.Range("R14").Value = 1489.981221456448
$ws.Range("S14").Value = 0.02493429112915455
$ws.Range("T14").Value = 0.02493429112915455

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Gnas"
$ws.Range("C15").Value = "Vipr1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 151.5583673333333
$ws.Range("H15").Value = 454.675102
$ws.Range("I15").Value = 0.1591332164896969
$ws.Range("J15").Value = 0.1591332164896969
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.146504
$ws.Range("N15").Value = 0.439512
$ws.Range("O15").Value = 0.02101489865429428
$ws.Range("P15").Value = 0.02101489865429428
$ws.Range("Q15").Value = 22.20390704780267
$ws.Range("R15").Value = 199.835163430224
$ws.Range("S15").Value = 0.003344168417062851
$ws.Range("T15").Value = 0.003344168417062851

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Gnas"
$ws.Range("C16").Value = "Vipr1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 151.5583673333333
$ws.Range("H16").Value = 454.675102
$ws.Range("I16").Value = 0.1591332164896969
$ws.Range("J16").Value = 0.1591332164896969
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.199477333333333
$ws.Range("N16").Value = 9.598431999999999
$ws.Range("O16").Value = 0.4589409975612385
$ws.Range("P16").Value = 0.4589409975612386
$ws.Range("Q16").Value = 484.9075609600071
$ws.Range("R16").Value = 4364.168048640064
$ws.Range("S16").Value = 0.07303275712091004
$ws.Range("T16").Value = 0.07303275712091004

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Gnas"
$ws.Range("C17").Value = "Vipr1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 151.5583673333333
$ws.Range("H17").Value = 454.675102
$ws.Range("I17").Value = 0.1591332164896969
$ws.Range("J17").Value = 0.1591332164896969
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.533112333333334
$ws.Range("N17").Value = 7.599337000000001
$ws.Range("O17").Value = 0.3633559422605724
$ws.Range("P17").Value = 0.3633559422605724
$ws.Range("Q17").Value = 383.9143695119305
$ws.Range("R17").Value = 3455.229325607375
$ws.Range("S17").Value = 0.05782199982256948
$ws.Range("T17").Value = 0.05782199982256948
